# Edit "SoC DFT 2016.pptx" slide 3 ("Background - Anna Chang"):
#  1. Update the title text with Chingwen's full first name.
#  2. Shift/resize the body placeholder (it grows a bit to fit a new line).
#  3. Add a new "Patents : ..." bullet paragraph right after the
#     "Programming Languages" bullet, before the trailing github link line.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# --- 1) Title -------------------------------------------------------------
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Background - Chingwen Anna Chang"

# --- 2) Body placeholder position / size ----------------------------------
$body = $s.Shapes.Item(2)
$body.Top = 80.13585
$body.Height = 295.98425196850394

# --- 3) New "Patents" paragraph --------------------------------------------
$tr = $body.TextFrame.TextRange

# Paragraph 15 is "Programming Languages : Perl, Python, ..."; insert the
# new paragraph right after it (and therefore right before the hyperlink
# paragraph that used to follow it).
$progLangPara = $tr.Paragraphs(15, 1)

$label = "Patents : "
$detail = "US8446161 (May, 2013)  US7499519 (Mar, 2009)  US 7545666 (Jun. 2009)"
[void]$progLangPara.InsertAfter("`r" + $label + $detail)

# Re-fetch the freshly created paragraph (now #16) and bold/color the label.
$tr2 = $body.TextFrame.TextRange
$patentsPara = $tr2.Paragraphs(16, 1)
$labelRange = $patentsPara.Characters(1, $label.Length)
$labelRange.Font.Bold = -1
$labelRange.Font.Color.RGB = 8865052
